$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Date Started value for row 6 (table "data")
$ws.Range("A6").Value = 44531

# Update the PctCompleted value for row 6
$ws.Range("D6").Value = 50

# Move the active selection on the frozen pane back to A2
$ws.Range("A2").Select()
